$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '43.669.88'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.332.09'
$ws.Range("E3").Value = '  +4.06%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '271.13'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '95.31'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").Value = '45.29'
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").Value = '8.14'
$ws.Range("E12").Value = '  +1.58%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = '2.679.90'
$ws.Range("E14").Value = '  +4.05%  '
$ws.Range("D15").Value = '15.60'
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("D16").Value = '0.860'
$ws.Range("E16").Value = '  +7.83%  '
$ws.Range("D17").Value = '2.330.06'
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").Value = '43.661.92'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '0.0000109'
$ws.Range("E19").Value = '  +5.16%  '
$ws.Range("D20").Value = '6.42'
$ws.Range("E20").Value = '  +6.79%  '
$ws.Range("D21").Value = '72.48'
$ws.Range("E21").Value = '  +2.75%  '
$ws.Range("D22").Value = '239.64'
$ws.Range("E22").Value = '  +2.96%  '
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("D24").Value = '9.37'
$ws.Range("E24").Value = '  +6.98%  '
$ws.Range("D26").Value = '2.53'
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").Value = '11.43'
$ws.Range("E27").Value = '  +2.16%  '
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("E30").Value = '  +9.00%  '
$ws.Range("D31").Value = '38.22'
$ws.Range("E31").Value = '  -2.92%  '
$ws.Range("D32").Value = '172.62'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '0.0899'
$ws.Range("E33").Value = '  -2.63%  '
$ws.Range("D34").Value = '5.49'
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("E36").Value = '  +3.33%  '
$ws.Range("E37").Value = '  -2.99%  '
$ws.Range("D38").Value = '4.38'
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").Value = '2.38'
$ws.Range("E40").Value = '  +10.21%  '
$ws.Range("E41").Value = '  +9.54%  '
$ws.Range("E42").Value = '  +18.94%  '
$ws.Range("D43").Value = '12.09'
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("D44").Value = '9.16'
$ws.Range("E44").Value = '  +8.92%  '
$ws.Range("D45").Value = '62.16'
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").Value = '5.35'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("E47").Value = '  +4.49%  '
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  +1.61%  '
$ws.Range("D50").Value = '0.191'
$ws.Range("E50").Value = '  +17.76%  '
$ws.Range("D51").Value = '2.558.12'
$ws.Range("E51").Value = '  +4.07%  '
